$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("logs")

# Common values shared across the new rows (logs 005 - 011)
$preprocess = 'remove break line, remove multiple spaces, convert unicode to ascii, convert to lower, space after punctuation, trim "space" and ","'
$features   = "13 features: #term, #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), #digit/#ascii, %kwName, %kwAddress, %kwPhone, #max_digit_skip_0 >= 7, #max_digit_skip_0 = 0, first_character_ascii, first_character_digit, last_character_ascii, last_character_digit"
$model      = "Neuron Network"
$modelDetails = "2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 2000"
$templateFilter = "0 filters: "

# New rows to append: Time, RunningTime(s), Test_Accuracy, Val_Accuracy, (blank score)
$newRows = @(
    @{ Time = "20160415_171634"; RunningTime = 1890.06;   TestAcc = 1;               ValAcc = 0.957095709570957; Score = 0.476744186046512 },
    @{ Time = "20160415_174804"; RunningTime = 1872.73;   TestAcc = 1;               ValAcc = 0.957095709570957; Score = 0.546511627906977 },
    @{ Time = "20160415_181916"; RunningTime = 1740.505;  TestAcc = 0.999333333333333; ValAcc = 0.95049504950495; Score = 0.119047619047619 },
    @{ Time = "20160415_184817"; RunningTime = 1324.746;  TestAcc = 1;               ValAcc = 0.957095709570957; Score = 0.476744186046512 },
    @{ Time = "20160415_191022"; RunningTime = 1216.841;  TestAcc = 1;               ValAcc = 0.957095709570957; Score = 0.27906976744186 }
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.Time
    $ws.Cells.Item($r, 2).Value = $row.RunningTime
    $ws.Cells.Item($r, 3).Value = $preprocess
    $ws.Cells.Item($r, 4).Value = $features
    $ws.Cells.Item($r, 5).Value = $model
    $ws.Cells.Item($r, 6).Value = $modelDetails
    $ws.Cells.Item($r, 7).Value = $row.TestAcc
    $ws.Cells.Item($r, 8).Value = $row.ValAcc
    $ws.Cells.Item($r, 9).Value = $templateFilter
    $ws.Cells.Item($r, 10).Value = $row.Score
}
